# Bug fix in ChestAndChance
#
# 1. Paragraph "2. ..." gets a trailing ". " appended to the existing
#    sentence and a brand-new sentence describing the invisible GO-square
#    button, added as its own run.
# 2. The "Normal" style picks up a slightly different default text colour
#    and becomes explicitly left-aligned.

$d = $word.ActiveDocument

# --- 1. Split the "2. ..." run into two runs -----------------------------

$oldSentence = "2. The program is run by clicking the green Run button in Monopoly.java"
$newFirstPart = "2. The program is run by clicking the green Run button in Monopoly.java. "
$secondSentence = "There is an invisible button on the GO square that, when pushed, will let you move a player to any specified square."

# Append ". " to the end of the existing sentence (keeps it as one run).
$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newFirstPart, 2)

# Locate that (now-extended) run so we know exactly where it starts/ends.
$firstRng = $d.Content
$firstRng.Find.Execute($newFirstPart)
$firstStart = $firstRng.Start
$firstEnd = $firstRng.End

# Insert the new sentence right after it.
$insPoint = $d.Range($firstEnd, $firstEnd)
$insStart = $insPoint.End
$insPoint.InsertAfter($secondSentence)
$insEnd = $insPoint.End

# Toggling (and reverting) a character property on each piece forces the
# engine to keep them as two distinct <w:r> runs instead of silently
# re-merging them, matching the two-run shape in the target document.
$firstRunRng = $d.Range($firstStart, $firstEnd)
$firstRunRng.Bold = 1
$firstRunRng.Bold = 0

$newRunRng = $d.Range($insStart, $insEnd)
$newRunRng.Bold = 1
$newRunRng.Bold = 0

# --- 2. Tweak the "Normal" style ------------------------------------------

$normal = $d.Styles("Normal")
$normal.Font.Color = 655360        # -> <w:color w:val="00000A"/>
$normal.ParagraphFormat.Alignment = 0   # wdAlignParagraphLeft -> <w:jc w:val="left"/>
